$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 84; this shifts the existing rows
# 84..164 down to 85..165 (including all of their cell values/styles).
$ws.Rows(84).Insert()

# Populate the newly inserted row 84 with its data.
$ws.Range("A84").Value = 11
$ws.Range("B84").Value = "Vega Monumental Concepción"
$ws.Range("C84").Value = "Bíobío"
$ws.Range("D84").Value = 44874
$ws.Range("E84").Value = 8
$ws.Range("F84").Value = 100112043
$ws.Range("G84").Value = "Pepino ensalada"
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 200
$ws.Range("K84").Value = 15000
$ws.Range("L84").Value = 16000
$ws.Range("M84").Value = 15500
$ws.Range("N84").Value = "$/caja 60 unidades"
$ws.Range("O84").Value = "Región de Arica y Parinacota"
$ws.Range("P84").Value = 258
$ws.Range("Q84").Value = 60
$ws.Range("R84").Value = "Hortaliza"
